# Adds the "zad 6.3" exercise block below the existing "zad 6.2" block,
# rewires the hidden Solver defined names to point at the new ranges, and
# lightly restyles the border of the tail of the "zad 6.2" constraint table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) New "zad 6.3" block: labels, decision vars, objective function
# ---------------------------------------------------------------------
$ws.Range("A23").Value = "zad 6.3"

$ws.Range("A24").Value = "zmienne decyzyjne:"
$ws.Range("D24").Value = "funkcja celu:"
$ws.Range("I24").Value = "zbiór dopuszczalny:"

$ws.Range("A25").Value = "x"
$ws.Range("B25").Value = "y"
$ws.Range("D25").Value = "f(x,y)=ax+by"

$ws.Range("I15").Copy()
$ws.Range("I25:M25").PasteSpecial(-4104)
$ws.Range("I25").Value = "przy x "
$ws.Range("J25").Value = "przy y"
$ws.Range("K25").Value = "znak"
$ws.Range("L25").Value = "ograniczenie"
$ws.Range("M25").Value = "formuła"

# Row 26: decision-variable values (yellow, boxed) + objective formula
$ws.Range("A5:B5").Copy()
$ws.Range("A26").PasteSpecial(-4104)
$ws.Range("A26").Value = 0.39999999999999997
$ws.Range("B26").Value = 0.60000000000000009

$ws.Range("D26").Value = "f(x,y)="

$ws.Range("E16").Copy()
$ws.Range("E26").PasteSpecial(-4104)
$ws.Range("E26").Formula = "=SUMPRODUCT(A26:B26,A28:B28)"

$ws.Range("H26").Value = "lepkość"
$ws.Range("I26").Value = 400
$ws.Range("J26").Value = 100
$ws.Range("K26").Value = ">="
$ws.Range("L26").Value = 200
$ws.Range("M26").Formula = "=SUMPRODUCT(`$A`$26:`$B`$26,I26:J26)"

# Row 27: "a)" style label + second objective line + constraint row 2
$ws.Range("A6").Copy()
$ws.Range("A27").PasteSpecial(-4104)
$ws.Range("A27").Value = "a"

$ws.Range("B27").Value = "b"

$ws.Range("D17").Copy()
$ws.Range("D27").PasteSpecial(-4104)
$ws.Range("D27").Value = "min"

$ws.Range("H27").Value = "Y(połysk)"
$ws.Range("I27").Value = 20
$ws.Range("J27").Value = 10
$ws.Range("K27").Value = ">="
$ws.Range("L27").Value = 14
$ws.Range("M27:M31").Formula = "=SUMPRODUCT(`$A`$26:`$B`$26,I27:J27)"

# Row 28: coefficients a,b + constraint row 3
$ws.Range("A7:B7").Copy()
$ws.Range("A28").PasteSpecial(-4104)
$ws.Range("A28").Value = 6
$ws.Range("B28").Value = 4

$ws.Range("H28").Value = "Z(trwałość)"
$ws.Range("I28").Value = 20
$ws.Range("J28").Value = 60
$ws.Range("K28").Value = ">="
$ws.Range("L28").Value = 30

# Row 29: constraint row 4 (equality x+y=1)
$ws.Range("H29").Value = "x+y=1"
$ws.Range("I29").Value = 1
$ws.Range("J29").Value = 1
$ws.Range("K29").Value = "="
$ws.Range("L29").Value = 1

# Row 30: constraint row 5
$ws.Range("I30").Value = 1
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = ">="
$ws.Range("L30").Value = 0

# Row 31: constraint row 6
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 1
$ws.Range("K31").Value = ">="
$ws.Range("L31").Value = 0

# ---------------------------------------------------------------------
# 2) Trailing notes below the new table
# ---------------------------------------------------------------------
$ws.Range("A33").Value = "b)"
$ws.Range("A34").Value = "zad 6.2"

# ---------------------------------------------------------------------
# 3) Restyle the tail of the "zad 6.2" box (open the bottom border up)
# ---------------------------------------------------------------------
$topBottom = $ws.Range("G20:K20")
$topBottom.Borders.Item(7).LineStyle = 1
$topBottom.Borders.Item(10).LineStyle = 1
$topBottom.Borders.Item(8).LineStyle = 1
$topBottom.Borders.Item(9).LineStyle = -4142

$topOnly = $ws.Range("G21:K21")
$topOnly.Borders.Item(7).LineStyle = -4142
$topOnly.Borders.Item(10).LineStyle = -4142
$topOnly.Borders.Item(9).LineStyle = -4142
$topOnly.Borders.Item(8).LineStyle = 1

# ---------------------------------------------------------------------
# 4) Column widths for the new columns H and L
# ---------------------------------------------------------------------
$ws.Columns.Item(8).ColumnWidth = 10.85546875
$ws.Columns.Item(12).ColumnWidth = 12.28515625

# ---------------------------------------------------------------------
# 5) View: selection + scroll position
# ---------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 10
$ws.Range("O29").Select()

# ---------------------------------------------------------------------
# 6) Rewire the hidden Solver parameters to the new "zad 6.3" ranges
# ---------------------------------------------------------------------
foreach ($n in $wb.Names) {
    switch ($n.Name) {
        "Arkusz1!solver_adj"  { $n.RefersTo = "=Arkusz1!`$A`$26:`$B`$26" }
        "Arkusz1!solver_lhs1" { $n.RefersTo = "=Arkusz1!`$M`$26" }
        "Arkusz1!solver_lhs2" { $n.RefersTo = "=Arkusz1!`$M`$27" }
        "Arkusz1!solver_lhs3" { $n.RefersTo = "=Arkusz1!`$M`$28" }
        "Arkusz1!solver_lhs4" { $n.RefersTo = "=Arkusz1!`$M`$29" }
        "Arkusz1!solver_num"  { $n.RefersTo = "=4" }
        "Arkusz1!solver_opt"  { $n.RefersTo = "=Arkusz1!`$E`$26" }
        "Arkusz1!solver_rel4" { $n.RefersTo = "=2" }
        "Arkusz1!solver_rhs1" { $n.RefersTo = "=Arkusz1!`$L`$26" }
        "Arkusz1!solver_rhs2" { $n.RefersTo = "=Arkusz1!`$L`$27" }
        "Arkusz1!solver_rhs3" { $n.RefersTo = "=Arkusz1!`$L`$28" }
        "Arkusz1!solver_rhs4" { $n.RefersTo = "=Arkusz1!`$L`$29" }
    }
}
